# "Generate Report for handoff"
#
# Row 3 of every sheet tracks b.md.md. A new handoff/handback round-trip
# happened for it, so:
#   - Status flips from "Handed back: in sync with en-US" to
#     "Ready for handoff" (Overview!B3/C3, zh-cn!B3, de-de!B3)
#   - The per-locale "Latest Handoff File" + "Latest Handoff Datetime"
#     cells (C3/D3 on the zh-cn and de-de sheets) now point at the new
#     handoff package for b.md.md instead of the stale a.md.md one.
#   - The matching hyperlink's display text is updated to match (the
#     link target itself is untouched).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet -----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-26 11:48:49"

foreach ($h in $zhcn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$C$3') {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# ---- de-de sheet -----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-26 11:49:00"

foreach ($h in $dede.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$C$3') {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
